# [WT A] Presentazione definitiva.
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 22 : "Tracciabilita' dei Design Goals" table (shape "Tabella 4")
# ---------------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$tblShape = $s22.Shapes.Item(2)
$tbl = $tblShape.Table

# Drop the last column ("CRITERI DELL'UTENTE FINALE" + the "/" filler cells
# below it) - the whole column is removed from the table.
$tbl.Columns.Item(5).Delete()

# Resize the remaining four columns to their final widths (EMU -> points).
$tbl.Columns.Item(1).Width = 2086749 / 12700
$tbl.Columns.Item(2).Width = 1765679 / 12700
$tbl.Columns.Item(3).Width = 1926214 / 12700
$tbl.Columns.Item(4).Width = 1926214 / 12700

# Row 2 ("ARCHITETTURA DEL SISTEMA ATTUALE") gains a new lead-in phrase.
$tbl.Cell(2,1).Shape.TextFrame.TextRange.Text = "DEFINIZIONE E IMPLEMENTAZIONE ARCHITETTURA DEL SISTEMA ATTUALE"

# Row 2 / column 2 ("  /" placeholder) now holds real analysis text.
$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text = "L’implementazione dei processi compiuti da genitori e personale soddisfa gli obiettivi in termini di tempi di risposta."

# Row 2 / column 3 ("            /" placeholder) now holds real analysis
# text, left-aligned.
$cell23 = $tbl.Cell(2,3).Shape.TextFrame.TextRange
$cell23.Text = "I controlli sull’input al’atto dell’inserimento (allo scopo di evitar failures) soddisfano gli obiettivi di ffidabilità e disponibilità."
$cell23.ParagraphFormat.Alignment = 1

# Row 3 / column 4 : "     /" -> " /"
$tbl.Cell(3,4).Shape.TextFrame.TextRange.Text = " /"

# ---------------------------------------------------------------------------
# Slide 23 : "SDD / Pregi e Difetti" recap slide
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)

$good = $s23.Shapes.Item(2)
$goodPara = $good.TextFrame.TextRange.Paragraphs(1)
$goodPara.Font.Bold = 1
$goodPara.Font.Italic = 1

$bad = $s23.Shapes.Item(3)
$badPara = $bad.TextFrame.TextRange.Paragraphs(1)
$badPara.Font.Bold = 1
$badPara.Font.Italic = 1
